# "added callsign when creating flight"
# Insert a new "CALLSIGN" row under the AIRCRAFT row on the PLAN_INFO sheet,
# update the AIRCRAFT registration, the DATE, and the TOF fuel value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLAN_INFO")

# Update the aircraft registration (was 5YSLN).
$ws.Range("B10").Value = "5YSLI"

# Insert a new row before current row 11 (REFUEL STOPS), shifting
# REFUEL STOPS / CREW / TOF rows down by one.
$ws.Rows.Item(11).Insert()

# New CALLSIGN row (label in column A mirrors the bold style used by the
# other field-label cells in column A, i.e. the style of A10/A12).
$ws.Range("A11").Value = "CALLSIGN"
$ws.Range("A11").Font.Bold = $true
$ws.Range("B11").Value = "SLI"

# Update the departure date (was 2024-02-02, serial 45324 -> serial 45480,
# i.e. 2024-07-07). Use Value2 with the raw serial number so the cell's
# existing custom date number format/style is left untouched (assigning a
# .NET DateTime to .Value can make Excel apply its own default date format).
$ws.Range("B8").Value2 = 45480

# Update TOF fuel at departure (now row 14 after the insert, was 1800).
$ws.Range("B14").Value = 1400

$ws.Range("B15").Select()
